$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell $ws "D2" "68.171.54"
Set-TextCell $ws "E2" "  -0.78%  "
Set-TextCell $ws "D3" "3.264.20"
Set-TextCell $ws "E3" "  -0.51%  "
Set-TextCell $ws "E4" "  +0.03%  "
Set-TextCell $ws "D5" "582.84"
Set-TextCell $ws "E5" "  -0.26%  "
Set-TextCell $ws "D6" "183.77"
Set-TextCell $ws "E6" "  -1.39%  "
Set-TextCell $ws "D7" "1.00"
Set-TextCell $ws "E7" "  +0.02%  "
Set-TextCell $ws "D8" "0.601"
Set-TextCell $ws "E8" "  +0.09%  "
Set-TextCell $ws "E9" "  -3.29%  "
Set-TextCell $ws "D10" "6.62"
Set-TextCell $ws "E10" "  -0.79%  "
Set-TextCell $ws "E11" "  -3.09%  "
Set-TextCell $ws "D12" "3.832.29"
Set-TextCell $ws "E12" "  -0.56%  "
Set-TextCell $ws "E13" "  +1.35%  "
Set-TextCell $ws "D14" "68.165.77"
Set-TextCell $ws "E14" "  -0.82%  "
Set-TextCell $ws "D15" "27.21"
Set-TextCell $ws "E15" "  -4.10%  "
Set-TextCell $ws "E16" "  -2.71%  "
Set-TextCell $ws "D17" "3.284.55"
Set-TextCell $ws "E17" "  +0.09%  "
Set-TextCell $ws "D18" "5.71"
Set-TextCell $ws "E18" "  -2.82%  "
Set-TextCell $ws "D19" "13.26"
Set-TextCell $ws "E19" "  -2.94%  "
Set-TextCell $ws "D20" "417.08"
Set-TextCell $ws "E20" "  +5.39%  "
Set-TextCell $ws "D21" "7.50"
Set-TextCell $ws "E21" "  -2.84%  "
Set-TextCell $ws "E22" "  +0.13%  "
Set-TextCell $ws "D23" "71.22"
Set-TextCell $ws "E23" "  -0.68%  "
Set-TextCell $ws "D24" "0.507"
Set-TextCell $ws "E24" "  -2.43%  "
Set-TextCell $ws "E25" "  -3.81%  "
Set-TextCell $ws "D26" "0.187"
Set-TextCell $ws "E26" "  -1.63%  "
Set-TextCell $ws "D27" "9.28"
Set-TextCell $ws "E27" "  -4.87%  "
Set-TextCell $ws "E28" "  +0.43%  "
Set-TextCell $ws "D29" "1.94"
Set-TextCell $ws "E29" "  -2.26%  "
Set-TextCell $ws "D30" "22.58"
Set-TextCell $ws "E30" "  -2.22%  "
Set-TextCell $ws "E31" "  -5.71%  "
Set-TextCell $ws "D32" "6.82"
Set-TextCell $ws "E32" "  -4.71%  "
Set-TextCell $ws "E33" "  -5.19%  "
Set-TextCell $ws "D34" "164.21"
Set-TextCell $ws "E34" "  +0.45%  "
Set-TextCell $ws "E35" "  -5.60%  "
Set-TextCell $ws "D36" "1.88"
Set-TextCell $ws "E36" "  -6.47%  "
Set-TextCell $ws "D37" "26.59"
Set-TextCell $ws "E37" "  -0.58%  "
Set-TextCell $ws "D38" "0.790"
Set-TextCell $ws "E38" "  -4.35%  "
Set-TextCell $ws "E39" "  -4.06%  "
Set-TextCell $ws "D40" "6.31"
Set-TextCell $ws "E40" "  -4.43%  "
Set-TextCell $ws "D41" "2.624.65"
Set-TextCell $ws "E41" "  -1.13%  "
Set-TextCell $ws "D42" "0.0672"
Set-TextCell $ws "E42" "  -2.83%  "
Set-TextCell $ws "E43" "  -5.54%  "
Set-TextCell $ws "D44" "334.98"
Set-TextCell $ws "E44" "  -1.97%  "
Set-TextCell $ws "D45" "24.06"
Set-TextCell $ws "E45" "  -6.23%  "
Set-TextCell $ws "E46" "  -3.82%  "
Set-TextCell $ws "D47" "6.22"
Set-TextCell $ws "E47" "  -2.34%  "
Set-TextCell $ws "D48" "0.978"
Set-TextCell $ws "E48" "  -1.54%  "
Set-TextCell $ws "D49" "0.1000"
Set-TextCell $ws "E49" "  -2.13%  "
Set-TextCell $ws "D50" "1.00"
Set-TextCell $ws "E50" "  +0.11%  "
Set-TextCell $ws "D51" "30.42"
Set-TextCell $ws "E51" "  -4.94%  "
